$d = $word.ActiveDocument

# Step 1: change the existing run's text from "1" to "2".
$null = $d.Content.Find.Execute("1", $true, $false, $false, $false, $false, $true, 1, $false, "2", 2)

# Step 2: append a brand-new run containing "1" right after the
# bookmarkEnd, while leaving the paragraph's own properties/attributes
# untouched. We do this by replacing the paragraph's content range with
# the fully-specified OOXML for the paragraph (same <w:p> attributes and
# <w:pPr> as before, the already-updated "2" run, the existing bookmark,
# and the new "1" run appended after the bookmark).
$r = $d.Content
$null = $r.InsertXML(@'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p w:rsidR="00D202AF" w:rsidRPr="003F5099" w:rsidRDefault="003F5099">
            <w:pPr>
              <w:rPr>
                <w:lang w:val="en-US"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:lang w:val="en-US"/>
              </w:rPr>
              <w:t>2</w:t>
            </w:r>
            <w:bookmarkStart w:id="0" w:name="_GoBack"/>
            <w:bookmarkEnd w:id="0"/>
            <w:r>
              <w:rPr>
                <w:lang w:val="en-US"/>
              </w:rPr>
              <w:t>1</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@)

Write-Host "Edit applied"
